# "Generate Report for Handback" - fills in the per-file localization
# handback columns (Latest Target File / Latest Handback File / Latest
# Handback DateTime) for both language sheets, flips the Overview status
# text, widens a few columns to fit the newly-populated long values, and
# adds hyperlinks on the new "Latest Target File" cells.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$targetMd1 = "31ec876b-71e2-48bc-8e51-ec81854b419e.md"
$targetMd2 = "b6188037-a32e-4f8f-8a32-1d79fa92826f.md"

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c0f134a7869c2ca824c4e457233fd6b1691986b9/e2e/"
$link1 = $ghBase + $targetMd1
$link2 = $ghBase + $targetMd2

# ---------------------------------------------------------------------
# Overview sheet: flip the per-language status text + widen the two
# status columns (E, F) to fit it.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
for ($r = 2; $r -le 3; $r++) {
    if ($wsOverview.Cells.Item($r, 5).Value() -eq $statusOld) {
        $wsOverview.Cells.Item($r, 5).Value = $statusNew
    }
    if ($wsOverview.Cells.Item($r, 6).Value() -eq $statusOld) {
        $wsOverview.Cells.Item($r, 6).Value = $statusNew
    }
}
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# Per-language sheets: zh-cn and de-de.  Each has the same layout:
#   A Source File Name   B File Extension   C Status
#   D Source Path        E Priority         F Content Duplicate
#   G Latest Handoff File  H Latest Handoff Datetime
#   I Latest Target File   J Latest Handback File   K Latest Handback DateTime
# Row 2 -> 31ec876b..., Row 3 -> b6188037...
# ---------------------------------------------------------------------
$langs = @(
    @{ Name = "zh-cn"; Ext = "zh-cn"; HandbackTime = "2016-09-07 11:06:49" },
    @{ Name = "de-de"; Ext = "de-de"; HandbackTime = "2016-09-07 11:06:57" }
)

$xlfHash1 = "60371152364574b67bb372ec43fdb92c6147132e"
$xlfHash2 = "8bc8d6b04852ec9338a4207c75a7fbec9f207ea3"

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Status column (C) text.
    for ($r = 2; $r -le 3; $r++) {
        if ($ws.Cells.Item($r, 3).Value() -eq $statusOld) {
            $ws.Cells.Item($r, 3).Value = $statusNew
        }
    }

    # Latest Target File (I) / Latest Handback File (J) / Latest Handback
    # DateTime (K) for both rows.
    $xlf1 = $targetMd1.Substring(0, $targetMd1.Length - 3) + "." + $xlfHash1 + "." + $lang.Ext + ".xlf"
    $xlf2 = $targetMd2.Substring(0, $targetMd2.Length - 3) + "." + $xlfHash2 + "." + $lang.Ext + ".xlf"

    $ws.Cells.Item(2, 9).Value = $targetMd1
    $ws.Cells.Item(2, 10).Value = $xlf1
    $ws.Cells.Item(2, 11).Value = $lang.HandbackTime

    $ws.Cells.Item(3, 9).Value = $targetMd2
    $ws.Cells.Item(3, 10).Value = $xlf2
    $ws.Cells.Item(3, 11).Value = $lang.HandbackTime

    # Widen the columns that now hold long file names (C, I, J).
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40

    # Rebuild the hyperlinks collection in row order (A2, I2, A3, I3) so
    # the relationship ids line up the way the handback report writes
    # them: A2=rId2 (unchanged), I2=rId3 (new), A3=rId4, I3=rId5 (new).
    $ws.Range("A1").Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $link1, "", "", $targetMd1)
    $ws.Hyperlinks.Add($ws.Range("I2"), $link1, "", "", $targetMd1)
    $ws.Hyperlinks.Add($ws.Range("A3"), $link2, "", "", $targetMd2)
    $ws.Hyperlinks.Add($ws.Range("I3"), $link2, "", "", $targetMd2)

    # Hyperlinks.Add stamps the built-in "Hyperlink" style on the cell;
    # restore the workbook's original custom HyperLink look (underlined,
    # cornflower blue, i.e. RGB 64,149,237 / #6495ED) on all four linked
    # cells -- Font.Color takes a BGR-packed value, hence 0xED9564.
    foreach ($cellRef in @("A2", "I2", "A3", "I3")) {
        $ws.Range($cellRef).Font.Underline = 2
        $ws.Range($cellRef).Font.Color = 0xED9564
    }
}
